# Add a "NewTag" (column D) WordNet-style POS grouping for every Penn
# Treebank tag row: adjectives -> "a", nouns -> "n", adverbs -> "r",
# verbs -> "v". Column D / header ("NewTag") already exists in the sheet;
# this only fills in the per-row values that were missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> grouping letter, keyed by the tag already present in column B
$groups = @{
    8  = "a"   # JJ   - Adjective
    9  = "a"   # JJR  - Adjective, comparative
    10 = "a"   # JJS  - Adjective, superlative
    13 = "n"   # NN   - Noun, singular or mass
    14 = "n"   # NNS  - Noun, plural
    15 = "n"   # NNP  - Proper noun, singular
    21 = "r"   # RB   - Adverb
    22 = "r"   # RBR  - Adverb, comparative
    23 = "r"   # RBS  - Adverb, superlative
    28 = "v"   # VB   - Verb, base form
    29 = "v"   # VBD  - Verb, past tense
    30 = "v"   # VBG  - Verb, gerund or present participle
    31 = "v"   # VBN  - Verb, past participle
    32 = "v"   # VBP  - Verb, non-3rd person singular present
    33 = "v"   # VBZ  - Verb, 3rd person singular present
}

foreach ($row in ($groups.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 4).Value = $groups[$row]
}

# Move the selection to where the author left off after the edit.
$ws.Range("D33").Select() | Out-Null
